$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9664817452430725
$ws.Range("B1").Value = 1.175292015075684
$ws.Range("C1").Value = 8.895249366760254
$ws.Range("D1").Value = 2.379584074020386
$ws.Range("E1").Value = 1.277235627174377
